$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 560, shifting existing rows 560-621 down to 561-622.
$ws.Rows.Item(560).Insert()

# Populate the newly inserted row 560 with the new data record.
$ws.Range("A560").Value = 9
$ws.Range("B560").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C560").Value = "Metropolitana"
$ws.Range("D560").Value = 45194
$ws.Range("E560").Value = 13
$ws.Range("F560").Value = 100112052
$ws.Range("G560").Value = "Albahaca"
$ws.Range("H560").Value = "Sin especificar"
$ws.Range("I560").Value = "Primera"
$ws.Range("J560").Value = 340
$ws.Range("K560").Value = 5000
$ws.Range("L560").Value = 5000
$ws.Range("M560").Value = 5000
$ws.Range("N560").Value = "$/paquete"
$ws.Range("O560").Value = "Región de Arica y Parinacota"
$ws.Range("P560").Value = 5000
$ws.Range("Q560").Value = 1
$ws.Range("R560").Value = "Hortaliza"
